# Trade #20 closed at 2026-02-17 13:18:12 - unknown UNKNOWN +0.000%
#
# Updates the Summary / Strategy Status roll-up numbers for the
# MarketMaking strategy and appends the newly-closed trade #20 row to
# both the "All Trades" and "MarketMaking" trade logs.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Summary sheet roll-up metrics
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.28   # Current Capital
$summary.Range("B4").Value = -0.73     # Total P&L $
$summary.Range("B5").Value = -0.73     # Total P&L %
$summary.Range("B6").Value = 20        # Total Trades
$summary.Range("B7").Value = 7         # Winning Trades
$summary.Range("B9").Value = 35        # Win Rate %

# ---------------------------------------------------------------
# 2. Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.28
$status.Range("D4").Value = 20
$status.Range("E4").Value = -0.73
$status.Range("F4").Value = -0.72
$status.Range("G4").Value = 35

# ---------------------------------------------------------------
# 3. Append new trade #20 row (row 21) to a trade-log sheet
# ---------------------------------------------------------------
function Add-TradeRow($ws) {
    $row = 21

    $ws.Cells.Item($row, 1).Value = 20

    # Date/time columns are stored as plain text in this workbook. A bare
    # "2026-02-17" would be auto-recognized and converted into a date
    # serial number, so it is entered with a leading apostrophe (just
    # like typing it into Excel) to force it to stay literal text.
    $ws.Cells.Item($row, 2).Value = "'2026-02-17"

    $ws.Cells.Item($row, 3).Value = "13:18:05"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "UP"
    $ws.Cells.Item($row, 6).Value = 0.25
    $ws.Cells.Item($row, 7).Value = 0.284314
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = 13.7255
    $ws.Cells.Item($row, 10).Value = 0.03
    $ws.Cells.Item($row, 11).Value = 99.28
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.14
}

Add-TradeRow($wb.Worksheets.Item("All Trades"))
Add-TradeRow($wb.Worksheets.Item("MarketMaking"))
